$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Invoices to collect": invoice rows were refreshed -- the old
# I-VS004002270/71 pair (dated 07 July 2018) is replaced with a new
# I-VS004002287/88/91 trio, and a new row is appended for the third one.
# ---------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("Invoices to collect")
$wsReturns = $wb.Worksheets.Item("Returns")

# Row 2 just gets a refreshed invoice number; formatting is untouched.
$wsInv.Range("A2").Value = "I-VS004002287"

# Row 3 used to hold a raw numeric collection amount (25000); it now
# matches the "Full" pattern used by row 2, so copy that row's look
# first and then overwrite the two values.
$wsInv.Range("A2:B2").Copy() | Out-Null
$wsInv.Range("A3").PasteSpecial(-4122) | Out-Null
$wsInv.Range("A3").Value = "I-VS004002288"
$wsInv.Range("B3").Value = "Full"

# Row 4 is brand new: pre-format A4 / B4 now (invoice-number look
# borrowed from A2/A3, numeric look borrowed from a numeric cell
# elsewhere in the workbook, Returns!A2); the values themselves are
# filled in below, after the "Payment methods" sheet updates.
$wsInv.Range("A4").PasteSpecial(-4122) | Out-Null

$wsReturns.Range("A2").Copy() | Out-Null
$wsInv.Range("B4").PasteSpecial(-4122) | Out-Null
$wsInv.Range("B4").Value = 26675

# ---------------------------------------------------------------------
# Sheet "Payment methods": correct the collection amount typo and push
# the cheque date out a month.
# ---------------------------------------------------------------------
$wsPay = $wb.Worksheets.Item("Payment methods")
$wsPay.Range("B2").Value = "7842.08"

# F3 and F4 are stored with quote-prefixed date-like text; re-apply the
# leading apostrophe so the cells keep their existing (text) style
# instead of Excel re-deriving a plain date style for them.
$wsPay.Range("F4").Value = "'07 August 2018"
$wsPay.Range("F3").Value = "'23 July 2018"

$wsPay.Activate() | Out-Null
$wsPay.Range("F6").Select() | Out-Null

# Back to "Invoices to collect" row 4's invoice number.
$wsInv.Range("A4").Value = "I-VS004002291"

# Restore "Invoices to collect" as the active/selected sheet & cell,
# matching the original workbook's active tab.
$wsInv.Activate() | Out-Null
$wsInv.Range("A4").Select() | Out-Null
